# Generate Report for Handoff
# Inserts a new "ready for handoff" row (file 0c05a52b-0c83-4586-9214-fbdd36cbc2b1)
# above the previously-existing row (file 8d18e275-99a5-4481-83b7-9dcffb239eac) on
# every sheet: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$newGuid = "0c05a52b-0c83-4586-9214-fbdd36cbc2b1"
$oldGuid = "8d18e275-99a5-4481-83b7-9dcffb239eac"
$newHash = "5ecd15e3fdae5f448c36f40c571c9b41a8f26c74"
$oldHash = "38bef17d5ca5f7108f12f6ff593a277e5c1978d7"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Push the existing data row (row 2) down to row 3, duplicating styles.
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-30-19 08:30:23"

# Hyperlinks don't follow the row shift automatically - rebuild them all.
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/6338203bb1567c2a4fe6ccf46ca35e8195a531b8/e2e/$newGuid.md", "", "", "$newGuid.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/6338203bb1567c2a4fe6ccf46ca35e8195a531b8/e2e/$oldGuid.md", "", "", "$oldGuid.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "$newGuid.$newHash.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-19 08:30:20"
$ws.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"

$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/6338203bb1567c2a4fe6ccf46ca35e8195a531b8/e2e/$newGuid.md", "", "", "$newGuid.md")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/6338203bb1567c2a4fe6ccf46ca35e8195a531b8/e2e/$newGuid.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8cc6d4387fd9d857c42c70350a6cb9a546af6e3d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newGuid.$newHash.zh-cn.xlf", "", "", "$newGuid.$newHash.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/6338203bb1567c2a4fe6ccf46ca35e8195a531b8/e2e/$oldGuid.md", "", "", "$oldGuid.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/6338203bb1567c2a4fe6ccf46ca35e8195a531b8/e2e/$oldGuid.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8cc6d4387fd9d857c42c70350a6cb9a546af6e3d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid.$oldHash.zh-cn.xlf", "", "", "$oldGuid.$oldHash.zh-cn.xlf")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "$newGuid.$newHash.de-de.xlf"
$ws.Range("E2").Value = "2016-03-19 08:30:23"
$ws.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"

$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/6338203bb1567c2a4fe6ccf46ca35e8195a531b8/e2e/$newGuid.md", "", "", "$newGuid.md")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/6338203bb1567c2a4fe6ccf46ca35e8195a531b8/e2e/$newGuid.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ef97ba5abf4bda8caec7f070251c79dbe1f0dfa1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newGuid.$newHash.de-de.xlf", "", "", "$newGuid.$newHash.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/6338203bb1567c2a4fe6ccf46ca35e8195a531b8/e2e/$oldGuid.md", "", "", "$oldGuid.md")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/6338203bb1567c2a4fe6ccf46ca35e8195a531b8/e2e/$oldGuid.md", "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ef97ba5abf4bda8caec7f070251c79dbe1f0dfa1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid.$oldHash.de-de.xlf", "", "", "$oldGuid.$oldHash.de-de.xlf")
